$wb = $excel.ActiveWorkbook
$wsOpt = $wb.Worksheets.Item("optimization_parameters")
$wsThr = $wb.Worksheets.Item("threshold_b")

# Insert a new row for the "L_curve" parameter right after the model/
# production_function row, shifting the remaining parameter rows down by one.
$wsOpt.Rows.Item(9).Insert()

# New row 9: L_curve parameter with default value 0.
$wsOpt.Cells.Item(9, 1).Value = "L_curve"
$wsOpt.Cells.Item(9, 2).Value = 0

# Rename the "Model" parameter label to "production_function".
$wsOpt.Cells.Item(8, 1).Value = "production_function"

# Make optimization_parameters the active sheet/selection (it was
# threshold_b before), with A8 selected.
$wsOpt.Activate()
$wsOpt.Range("A8").Select()
